$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.049.91'
$ws.Range('E2').Value = '  -0.06%  '
$ws.Range('D3').Value = '1.788.45'
$ws.Range('E3').Value = '  -0.03%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '226.98'
$ws.Range('E5').Value = '  +1.81%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.544'
$ws.Range('E6').Value = '  -1.24%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '32.18'
$ws.Range('E8').Value = '  -0.35%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.295'
$ws.Range('E9').Value = '  +3.72%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0685'
$ws.Range('E10').Value = '  -3.70%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0941'
$ws.Range('E11').Value = '  +1.15%  '
$ws.Range('D12').Value = '2.046.61'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '11.30'
$ws.Range('E13').Value = '  +2.38%  '
$ws.Range('D14').Value = '1.798.55'
$ws.Range('E14').Value = '  +0.67%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.624'
$ws.Range('D16').Value = '34.048.70'
$ws.Range('E16').Value = '  -0.11%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '4.19'
$ws.Range('E17').Value = '  +0.50%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '67.87'
$ws.Range('E18').Value = '  -0.11%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '242.16'
$ws.Range('E19').Value = '  -0.82%  '
$ws.Range('E20').Value = '  -1.01%  '
$ws.Range('E21').Value = '  +0.07%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '10.72'
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('E23').Value = '  +0.52%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.06'
$ws.Range('E24').Value = '  -2.67%  '
$ws.Range('E25').Value = '  +1.83%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.14'
$ws.Range('E26').Value = '  +1.66%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '16.19'
$ws.Range('E27').Value = '  -0.58%  '
$ws.Range('E28').Value = '  +0.93%  '
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('E30').Value = '  +2.02%  '
$ws.Range('E31').Value = '  -0.60%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.64'
$ws.Range('E32').Value = '  -0.54%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.61'
$ws.Range('E33').Value = '  +3.66%  '
$ws.Range('E34').Value = '  +1.59%  '
$ws.Range('D35').Value = '1.402.89'
$ws.Range('E35').Value = '  +1.48%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.652'
$ws.Range('E36').Value = '  +0.51%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.05'
$ws.Range('E37').Value = '  -0.42%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0188'
$ws.Range('E38').Value = '  +1.87%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.33'
$ws.Range('E39').Value = '  +7.50%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '79.83'
$ws.Range('E40').Value = '  +0.08%  '
$ws.Range('E41').Value = '  +0.34%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.918'
$ws.Range('E42').Value = '  +0.37%  '
$ws.Range('E43').Value = '  -0.40%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '13.62'
$ws.Range('E44').Value = '  +13.64%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '6.13'
$ws.Range('E45').Value = '  +3.40%  '
$ws.Range('D46').Value = '0.0₆0140'
$ws.Range('E46').Value = '  +5.49%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0507'
$ws.Range('E47').Value = '  +1.59%  '
$ws.Range('E48').Value = '  +2.35%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '107.52'
$ws.Range('E49').Value = '  -0.07%  '
$ws.Range('D50').Value = '1.948.21'
$ws.Range('E50').Value = '  +0.13%  '
$ws.Range('E51').Value = '  +0.01%  '

Write-Host "Applied 87 cell updates"
